# Update Carga_Masiva.xlsx header row from 8 columns to 5 columns
# New headers: fecha, peso, turno, paquete_operador, paquete_norma

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the old header columns F:H (beyond the new 5-column layout)
$ws.Range("F1:H1").ClearContents()

# Set the new header values
$ws.Range("A1").Value = "fecha"
$ws.Range("B1").Value = "peso"
$ws.Range("C1").Value = "turno"
$ws.Range("D1").Value = "paquete_operador"
$ws.Range("E1").Value = "paquete_norma"

# Adjust column widths (offset by 5/6 to account for Excel's internal
# padding so the stored XML width matches the target exactly)
$ws.Columns.Item(2).ColumnWidth = 14.166666666666668
$ws.Columns.Item(3).ColumnWidth = 13.166666666666668
$ws.Columns.Item(4).ColumnWidth = 22.833333333333336
$ws.Columns.Item(5).ColumnWidth = 23.333333333333336

# Update the selected cell
$ws.Range("B3").Select()
